$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/31/2023  Through  8/6/2023"

# --- Weekly crime-stat table updates (rows 14-30) ---
# Row 14
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = -50
$ws.Range("J14").Value = 21
$ws.Range("K14").Value = -47.619047619047
# Row 15
$ws.Range("C15").Value = 2
$ws.Range("E15").Value = -60
$ws.Range("F15").Value = 10
$ws.Range("G15").Value = 19
$ws.Range("H15").Value = -47.368421052631
$ws.Range("I15").Value = 72
$ws.Range("J15").Value = 110
$ws.Range("K15").Value = -34.545454545454
$ws.Range("L15").Value = -25.773195876288
$ws.Range("M15").Value = 35.849056603773
$ws.Range("N15").Value = -47.826086956521
# Row 16
$ws.Range("C16").Value = 37
$ws.Range("D16").Value = 42
$ws.Range("E16").Value = -11.904761904761
$ws.Range("F16").Value = 146
$ws.Range("G16").Value = 164
$ws.Range("H16").Value = -10.975609756097
$ws.Range("I16").Value = 1059
$ws.Range("J16").Value = 1240
$ws.Range("K16").Value = -14.596774193548
$ws.Range("L16").Value = 26.372315035799
$ws.Range("M16").Value = 45.867768595041
$ws.Range("N16").Value = -82.996146435452
# Row 17
$ws.Range("C17").Value = 65
$ws.Range("D17").Value = 44
$ws.Range("E17").Value = 47.727272727272
$ws.Range("F17").Value = 206
$ws.Range("G17").Value = 179
$ws.Range("H17").Value = 15.083798882681
$ws.Range("I17").Value = 1323
$ws.Range("J17").Value = 1213
$ws.Range("K17").Value = 9.068425391591
$ws.Range("L17").Value = 19.189189189189
$ws.Range("M17").Value = 72.041612483745
$ws.Range("N17").Value = -35.651750972762
# Row 18
$ws.Range("C18").Value = 35
$ws.Range("D18").Value = 51
$ws.Range("E18").Value = -31.372549019607
$ws.Range("F18").Value = 149
$ws.Range("G18").Value = 239
$ws.Range("H18").Value = -37.65690376569
$ws.Range("I18").Value = 1249
$ws.Range("J18").Value = 1809
$ws.Range("K18").Value = -30.956329463792
$ws.Range("L18").Value = 10.726950354609
$ws.Range("M18").Value = 15.541165587419
$ws.Range("N18").Value = -82.323804132465
# Row 19
$ws.Range("C19").Value = 227
$ws.Range("E19").Value = -13.026819923371
$ws.Range("F19").Value = 896
$ws.Range("G19").Value = 1036
$ws.Range("H19").Value = -13.513513513513
$ws.Range("I19").Value = 6792
$ws.Range("J19").Value = 6743
$ws.Range("K19").Value = 0.726679519501
$ws.Range("L19").Value = 69.249937702467
$ws.Range("M19").Value = 10.456984875589
$ws.Range("N19").Value = -64.70037939816
# Row 20
$ws.Range("C20").Value = 13
$ws.Range("D20").Value = 26
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 60
$ws.Range("G20").Value = 81
$ws.Range("H20").Value = -25.925925925925
$ws.Range("I20").Value = 353
$ws.Range("J20").Value = 424
$ws.Range("K20").Value = -16.745283018867
$ws.Range("L20").Value = 4.747774480712
$ws.Range("M20").Value = 58.29596412556
$ws.Range("N20").Value = -90.078695896571
# Row 21
$ws.Range("C21").Value = 379
$ws.Range("D21").Value = 430
$ws.Range("E21").Value = -11.860465116279
$ws.Range("F21").Value = 1468
$ws.Range("G21").Value = 1720
$ws.Range("H21").Value = -14.651162790697
$ws.Range("I21").Value = 10859
$ws.Range("J21").Value = 11560
$ws.Range("K21").Value = -6.06401384083
$ws.Range("L21").Value = 44.056778986468
$ws.Range("M21").Value = 20.53502053502
$ws.Range("N21").Value = -71.674883272034
# Row 22
$ws.Range("C22").Value = 15
$ws.Range("D22").Value = 8
$ws.Range("E22").Value = 87.5
$ws.Range("F22").Value = 49
$ws.Range("G22").Value = 38
$ws.Range("H22").Value = 28.947368421052
$ws.Range("I22").Value = 389
$ws.Range("J22").Value = 404
$ws.Range("K22").Value = -3.712871287128
$ws.Range("L22").Value = 44.609665427509
$ws.Range("M22").Value = 23.885350318471
# Row 23
$ws.Range("C23").Value = 6
$ws.Range("E23").Value = -40
$ws.Range("F23").Value = 46
$ws.Range("G23").Value = 39
$ws.Range("H23").Value = 17.948717948717
$ws.Range("I23").Value = 245
$ws.Range("J23").Value = 274
$ws.Range("K23").Value = -10.583941605839
$ws.Range("L23").Value = -18.60465116279
$ws.Range("M23").Value = 9.865470852017
# Row 24
$ws.Range("C24").Value = 425
$ws.Range("D24").Value = 430
$ws.Range("E24").Value = -1.162790697674
$ws.Range("F24").Value = 1699
$ws.Range("G24").Value = 1824
$ws.Range("H24").Value = -6.853070175438
$ws.Range("I24").Value = 12175
$ws.Range("J24").Value = 12840
$ws.Range("K24").Value = -5.179127725856
$ws.Range("L24").Value = 54.074917742343
$ws.Range("M24").Value = 21.762176217621
# Row 25
$ws.Range("C25").Value = 91
$ws.Range("D25").Value = 95
$ws.Range("E25").Value = -4.210526315789
$ws.Range("F25").Value = 387
$ws.Range("G25").Value = 406
$ws.Range("H25").Value = -4.679802955665
$ws.Range("I25").Value = 2967
$ws.Range("J25").Value = 2767
$ws.Range("K25").Value = 7.228044813877
$ws.Range("L25").Value = 37.043879907621
$ws.Range("M25").Value = 38.192827200745
# Row 26
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 17
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = -29.166666666666
$ws.Range("I26").Value = 130
$ws.Range("J26").Value = 169
$ws.Range("K26").Value = -23.076923076923
$ws.Range("L26").Value = -9.090909090909
# Row 27
$ws.Range("C27").Value = 18
$ws.Range("D27").Value = 28
$ws.Range("E27").Value = -35.714285714285
$ws.Range("F27").Value = 69
$ws.Range("G27").Value = 89
$ws.Range("H27").Value = -22.471910112359
$ws.Range("I27").Value = 543
$ws.Range("J27").Value = 578
$ws.Range("K27").Value = -6.055363321799
$ws.Range("L27").Value = 19.340659340659
# Row 28
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 75
$ws.Range("I28").Value = 28
$ws.Range("K28").Value = -24.324324324324
$ws.Range("L28").Value = 3.703703703703
$ws.Range("M28").Value = 40
$ws.Range("N28").Value = -67.441860465116
# Row 29
$ws.Range("F29").Value = 5
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 25
$ws.Range("I29").Value = 23
$ws.Range("K29").Value = -28.125
$ws.Range("L29").Value = -8
$ws.Range("M29").Value = 43.75
$ws.Range("N29").Value = -69.736842105263
# Row 30
$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 5
$ws.Range("E30").Value = -80
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 19
$ws.Range("H30").Value = -89.473684210526
$ws.Range("I30").Value = 60
$ws.Range("J30").Value = 111
$ws.Range("K30").Value = -45.945945945945
$ws.Range("L30").Value = -36.170212765957
